$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty row 6 with the result row that was
# accidentally blanked out
$ws.Range("A6").Value = 6578
$ws.Range("B6").Value = "Krokker Mihály"
$ws.Range("C6").Value = "Pandúr Lövész-Klub Sportegyesület"

# Re-add the empty row underneath (row 7) that keeps the competition id
# (the number is stored as text, as it was before)
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "8653"
$ws.Range("A7").Style = "Normal"
$ws.Range("V7").Value = "VID_00001"
